$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the filtered_feature_bc_matrix.h5 path to cell_feature_matrix.h5
$ws.Range("E2").Value = "datasets/10x_xenium_mouse_pup_preview/cell_feature_matrix.h5"

# Update the kmeans clusters.csv path to the new clustering path
$ws.Range("G2").Value = "datasets/10x_xenium_mouse_pup_preview/analysis/clustering/gene_expression_kmeans_10_clusters/clusters.csv"

# Move the active selection to E3
$ws.Range("E3").Select()
